# Append one new observation record as row 5 of the "Artfynd" sheet,
# following the same column layout as the existing rows (1 header row +
# rows 2-4 of data). This grows the used range from A1:AY4 to A1:AY5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 5

$ws.Cells.Item($row, 1).Value  = 112379124              # A  Id
$ws.Cells.Item($row, 2).Value  = 89503                   # B  Taxonsorteringsordning
$ws.Cells.Item($row, 3).Value  = "Ovaliderad"             # C  Valideringsstatus
$ws.Cells.Item($row, 4).Value  = "LC"                     # D  Rödlistade
$ws.Cells.Item($row, 5).Value  = 5447                     # E  TaxonId
$ws.Cells.Item($row, 6).Value  = "Vedticka"                # F  Artnamn
$ws.Cells.Item($row, 7).Value  = "Fuscoporia viticola"      # G  Vetenskapligt namn
$ws.Cells.Item($row, 8).Value  = "(Schwein.) Murrill"        # H  Auktor
$ws.Cells.Item($row, 9).Value  = ""                       # I  Antal (blank)

$ws.Cells.Item($row, 16).Value = "Lortmossen, Vstm"        # P  Lokalnamn
$ws.Cells.Item($row, 17).Value = 531712                   # Q  Ost
$ws.Cells.Item($row, 18).Value = 6622531                  # R  Nord
$ws.Cells.Item($row, 19).Value = 25                        # S  Noggrannhet
$ws.Cells.Item($row, 20).Value = "Västmanland"              # T  Län
$ws.Cells.Item($row, 21).Value = "Skinnskatteberg"          # U  Kommun
$ws.Cells.Item($row, 22).Value = "Västmanland"              # V  Provins
$ws.Cells.Item($row, 23).Value = "Skinnskatteberg"          # W  Församling

# Y (Startdatum) / AA (Slutdatum) hold plain date-looking text, not real
# dates. Pre-format the cell as Text ("@") before writing so Excel does not
# silently convert the string into a date serial number, then restore the
# Normal style so the cell is left without a lingering number format.
$ws.Cells.Item($row, 25).NumberFormat = "@"
$ws.Cells.Item($row, 25).Value = "2023-08-27"              # Y  Startdatum
$ws.Cells.Item($row, 25).Style = "Normal"

$ws.Cells.Item($row, 27).NumberFormat = "@"
$ws.Cells.Item($row, 27).Value = "2023-09-18"              # AA Slutdatum
$ws.Cells.Item($row, 27).Style = "Normal"

$ws.Cells.Item($row, 30).Value = $false                    # AD Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false                    # AE Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false                    # AG Ospontan

$ws.Cells.Item($row, 46).Value = ""                        # AT Bestämningsår (blank)

$ws.Cells.Item($row, 49).Value = "Mikael Hagström"          # AW Rapportör
$ws.Cells.Item($row, 50).Value = "Mikael Hagström"          # AX Observatörer
$ws.Cells.Item($row, 51).Value = ""                        # AY Projektnamn (blank)
